# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" footer timestamp (08:05 -> 08:35) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 08:35"

# --- Estados Unidos (row 4) : refreshed totals ---
$ws.Range("B4").Value = 1621134
$ws.Range("C4").Value = 232
$ws.Range("D4").Value = 382202
$ws.Range("E4").Value = 1142575
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 96357

# --- Ucrania / Polonia swap places in the ranking (rows 35-36) ---
# Row 35 becomes Ucrania with its refreshed totals
$ws.Range("A35").Value = "Ucrania"
$ws.Range("B35").Value = 20148
$ws.Range("C35").Value = 442
$ws.Range("D35").Value = 6585
$ws.Range("E35").Value = 12975
$ws.Range("G35").Value = 9
$ws.Range("H35").Value = 588

# Row 36 becomes Polonia, keeping its previous totals
$ws.Range("A36").Value = "Polonia"
$ws.Range("B36").Value = 20143
$ws.Range("D36").Value = 8452
$ws.Range("E36").Value = 10719
$ws.Range("H36").Value = 972

# --- Georgia (row 122) : refreshed totals ---
$ws.Range("B122").Value = 723
$ws.Range("C122").Value = 2
$ws.Range("D122").Value = 495
$ws.Range("E122").Value = 216
